# "updated Kevin results tab"
#
# Updates the data in the "Kevin's Tests" worksheet (new timing figures for
# rows 2-5, plus a handful of revised raw counts in row 4), and moves the
# active/selected tab from "Greg's tests" over to "Kevin's Tests".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kevin's Tests")

# --- Row 2 (Levenshtein's Distance) : revised timing columns -----------
$ws.Range("L2").Value = 0.154304
$ws.Range("M2").Value = 0.164076
$ws.Range("N2").Value = 0.165051

# --- Row 3 (Jaro-Winkler Distance) : revised timing columns ------------
$ws.Range("L3").Value = 0.022464
$ws.Range("M3").Value = 0.0634762
$ws.Range("N3").Value = 0.0625023

# --- Row 4 (Hunt-McIlroy Distance) : revised raw counts + timings ------
$ws.Range("B4").Value = 84
$ws.Range("C4").Value = 55
$ws.Range("D4").Value = 55
$ws.Range("F4").Value = 55
$ws.Range("G4").Value = 55
$ws.Range("H4").Value = 88
$ws.Range("I4").Value = 97
$ws.Range("J4").Value = 52
$ws.Range("K4").Value = 52
$ws.Range("L4").Value = 0.0244218
$ws.Range("M4").Value = 0.06152930
$ws.Range("N4").Value = 0.036135

# --- Row 5 (Needleman-Wunsch Distance) : revised timing columns --------
$ws.Range("L5").Value = 0.112313
$ws.Range("M5").Value = 0.11622
$ws.Range("N5").Value = 0.108406

# --- Selection / active tab --------------------------------------------
# Kevin's results tab becomes the active sheet/selection; the range
# A1:O5 (the whole table) ends up selected, anchored at O5.
$ws.Activate() | Out-Null
$ws.Range("A1:O5").Select() | Out-Null
